# Refresh Carbuncle_Profits market-data snapshot (currentAveragePrice* / LevePrice* / LeveProfit* columns)
# across the eight crafter sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Updates H:N for the affected leve rows; some rows gain/lose a trailing M or N cell
# because the sign of the computed profit crossed zero (cleared with ClearContents()).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1037.8572
$ws.Range("I18").Value = 500
$ws.Range("J18").Value = 1441.25
$ws.Range("K18").Value = 500
$ws.Range("L18").Value = 1441.25
$ws.Range("M18").Value = -216
$ws.Range("N18").Value = -2009.25

$ws.Range("H129").Value = 1011.9804
$ws.Range("I129").Value = 421.125
$ws.Range("J129").Value = 1121.907
$ws.Range("K129").Value = 1263.375
$ws.Range("L129").Value = 3365.721
$ws.Range("M129").Value = 3736.625
$ws.Range("N129").Value = -13365.721

$ws.Range("H132").Value = 1414.0869
$ws.Range("I132").Value = 1131.8948
$ws.Range("J132").Value = 2754.5
$ws.Range("K132").Value = 3395.6844
$ws.Range("L132").Value = 8263.5
$ws.Range("M132").Value = -865.6844000000001
$ws.Range("N132").Value = -13323.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2487.375
$ws.Range("I2").Value = 2758.4285
$ws.Range("J2").Value = 590
$ws.Range("K2").Value = 2758.4285
$ws.Range("L2").Value = 590
$ws.Range("M2").Value = -2645.4285
$ws.Range("N2").Value = -816

$ws.Range("H32").Value = 4202.5615
$ws.Range("I32").Value = 2674.8088
$ws.Range("K32").Value = 2674.8088
$ws.Range("M32").Value = -2387.8088

$ws.Range("H45").Value = 932.4706
$ws.Range("I45").Value = 873.4666999999999
$ws.Range("J45").Value = 1375
$ws.Range("K45").Value = 873.4666999999999
$ws.Range("L45").Value = 1375
$ws.Range("M45").Value = -496.4666999999999
$ws.Range("N45").Value = -2129

$ws.Range("H74").Value = 1250.421
$ws.Range("I74").Value = 1287.8572
$ws.Range("J74").Value = 1145.6
$ws.Range("K74").Value = 1287.8572
$ws.Range("L74").Value = 1145.6
$ws.Range("M74").Value = -413.8571999999999
$ws.Range("N74").Value = -2893.6

$ws.Range("H77").Value = 1250.421
$ws.Range("I77").Value = 1287.8572
$ws.Range("J77").Value = 1145.6
$ws.Range("K77").Value = 6439.286
$ws.Range("L77").Value = 5728
$ws.Range("M77").Value = -2071.286
$ws.Range("N77").Value = -14464

$ws.Range("H97").Value = 3331.0908
$ws.Range("I97").Value = 2715
$ws.Range("K97").Value = 2715
$ws.Range("M97").Value = -2219

$ws.Range("H116").Value = 2487.375
$ws.Range("I116").Value = 2758.4285
$ws.Range("J116").Value = 590
$ws.Range("K116").Value = 2758.4285
$ws.Range("L116").Value = 590
$ws.Range("M116").Value = -464.4285
$ws.Range("N116").Value = -5178

$ws.Range("H122").Value = 2435.6316
$ws.Range("I122").Value = 2035.9286
$ws.Range("J122").Value = 3554.8
$ws.Range("K122").Value = 6107.7858
$ws.Range("L122").Value = 10664.4
$ws.Range("M122").Value = -3657.7858
$ws.Range("N122").Value = -15564.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2487.375
$ws.Range("I3").Value = 2758.4285
$ws.Range("J3").Value = 590
$ws.Range("K3").Value = 2758.4285
$ws.Range("L3").Value = 590
$ws.Range("M3").Value = -2644.4285
$ws.Range("N3").Value = -818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2596.0386
$ws.Range("I31").Value = 2341.5
$ws.Range("J31").Value = 2730.7942
$ws.Range("K31").Value = 2341.5
$ws.Range("L31").Value = 2730.7942
$ws.Range("M31").Value = -2046.5
$ws.Range("N31").Value = -3320.7942

$ws.Range("H34").Value = 2596.0386
$ws.Range("I34").Value = 2341.5
$ws.Range("J34").Value = 2730.7942
$ws.Range("K34").Value = 2341.5
$ws.Range("L34").Value = 2730.7942
$ws.Range("M34").Value = -2139.5
$ws.Range("N34").Value = -3134.7942

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H63").Value = 31000
$ws.Range("J63").Value = 31000
$ws.Range("L63").Value = 31000
$ws.Range("N63").Value = -32372

$ws.Range("H66").Value = 31000
$ws.Range("J66").Value = 31000
$ws.Range("L66").Value = 93000
$ws.Range("N66").Value = -99864

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1237.3334
$ws.Range("I64").Value = 1237.3334
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3712.0002
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -3442.0002
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 1237.3334
$ws.Range("I67").Value = 1237.3334
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3712.0002
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2776.0002
$ws.Range("N67").ClearContents()

$ws.Range("H129").Value = 1544.8286
$ws.Range("I129").Value = 768.82355
$ws.Range("J129").Value = 2277.7222
$ws.Range("K129").Value = 2306.47065
$ws.Range("L129").Value = 6833.1666
$ws.Range("M129").Value = 2693.52935
$ws.Range("N129").Value = -16833.1666

$ws.Range("H131").Value = 5060
$ws.Range("I131").Value = 384.2857
$ws.Range("J131").Value = 6878.3335
$ws.Range("K131").Value = 1152.8571
$ws.Range("L131").Value = 20635.0005
$ws.Range("M131").Value = 3887.1429
$ws.Range("N131").Value = -30715.0005

$ws.Range("H137").Value = 1668.919
$ws.Range("I137").Value = 1454.6666
$ws.Range("J137").Value = 1871.8948
$ws.Range("K137").Value = 4363.9998
$ws.Range("L137").Value = 5615.6844
$ws.Range("M137").Value = 736.0002000000004
$ws.Range("N137").Value = -15815.6844

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 474.4
$ws.Range("I9").Value = 474.4
$ws.Range("K9").Value = 474.4
$ws.Range("M9").Value = -250.4

$ws.Range("H22").Value = 1250
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 1328.5714
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 1328.5714
$ws.Range("M22").Value = -405
$ws.Range("N22").Value = -1918.5714

$ws.Range("H27").Value = 1250
$ws.Range("I27").Value = 700
$ws.Range("J27").Value = 1328.5714
$ws.Range("K27").Value = 700
$ws.Range("L27").Value = 1328.5714
$ws.Range("M27").Value = -593
$ws.Range("N27").Value = -1542.5714

$ws.Range("H98").Value = 31773
$ws.Range("J98").Value = 31773
$ws.Range("L98").Value = 31773
$ws.Range("N98").Value = -37763

$ws.Range("H139").Value = 37160
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 37160
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 37160
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -47440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 59979.332
$ws.Range("J135").Value = 59979.332
$ws.Range("L135").Value = 59979.332
$ws.Range("N135").Value = -70119.33199999999

$ws.Range("H137").Value = 74081.664
$ws.Range("J137").Value = 74081.664
$ws.Range("L137").Value = 74081.664
$ws.Range("N137").Value = -84281.664

$ws.Range("H139").Value = 72000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 72000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 72000
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -82280

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()
